$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: collapse the three runs ". " / "Run the application" / "."
# (immediately after the closing curly quote) into a single run reading
# ". Run the application." without touching the preceding quote run.
# ---------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(". Run the application.") | Out-Null
$start1 = $rng1.Start
$end1 = $rng1.End
$d.Range($start1, $end1).Delete()
$d.Range($start1, $start1).InsertAfter(". Run the application.")

# ---------------------------------------------------------------
# Change 2: collapse the two italic runs "Assignment_WebTech_MERNStack_
# WEB027_ReactJS" and ".docx" into one italic run, keeping the
# surrounding curly-quote runs untouched.
# ---------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Assignment_WebTech_MERNStack_WEB027_ReactJS.docx") | Out-Null
$start2 = $rng2.Start
$end2 = $rng2.End
$target2 = $d.Range($start2, $end2)
# Round-trip through a distinct value first so the engine treats this as a
# genuine text change (a same-value assignment is a no-op) and rewrites the
# whole selected span as one run, inheriting the (italic) formatting of the
# range; then restore the exact original text.
$target2.Text = "Assignment_WebTech_MERNStack_WEB027_ReactJS.docx#TMP#"
$target2.Text = "Assignment_WebTech_MERNStack_WEB027_ReactJS.docx"

# ---------------------------------------------------------------
# Change 3: "For this purpose make necessary changes..." becomes
# "For this " / "purpose," / " make necessary changes..." — i.e. a comma
# is added after "purpose" and the sentence ends up split across three
# runs at that point.
# ---------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("For this purpose make necessary changes in React JS application, which was previously requesting JSON server.") | Out-Null
$sentenceStart = $rng3.Start

$purposeStart = $sentenceStart + ("For this ").Length
$purposeEnd = $purposeStart + ("purpose").Length

# Insert the comma right after "purpose". Anchoring the insertion point
# with a bookmark (and inserting through the bookmark's own Range) forces
# a run boundary there, so "purpose" and the trailing text do not get
# silently re-merged by the engine's run-coalescing.
$d.Bookmarks.Add("commaAnchor", $d.Range($purposeEnd, $purposeEnd)) | Out-Null
$d.Bookmarks("commaAnchor").Range.InsertAfter(",")
$d.Bookmarks("commaAnchor").Delete()

# Then split "For this " away from "purpose," so each becomes its own run.
$d.Bookmarks.Add("leadAnchor", $d.Range($sentenceStart, $purposeStart)) | Out-Null
$d.Bookmarks("leadAnchor").Delete()
